$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update header text: [ug] -> [mg]
$ws.Range("B1").Value = "Glutamate par semaine [mg]"

# Update Glutamate values (column B, rows 2-79) per new nutrition data
$values = @{
    2 = 118439.825
    3 = 124744.43
    4 = 68241.38
    5 = 53458.72499999999
    6 = 110010.68
    7 = 86093.68000000001
    8 = 132777.345
    9 = 100863.22
    10 = 117148.05
    11 = 166689.29
    12 = 123913.49
    13 = 90563.24500000001
    14 = 85980.2
    15 = 65107.255
    16 = 120467.95
    17 = 158445.38
    18 = 154194.18
    19 = 226822.09
    20 = 100577.305
    21 = 53085.46000000001
    22 = 95929.57000000001
    23 = 59157.375
    24 = 124697.145
    25 = 98162.85000000001
    26 = 73239.39999999999
    27 = 72798.32500000001
    28 = 90364.04999999999
    29 = 104696.46
    30 = 191426.05
    31 = 116659.655
    32 = 57688.325
    33 = 130326.05
    34 = 63057.87499999999
    35 = 76093.2
    36 = 59518.98
    37 = 157793.455
    38 = 290564.695
    39 = 109355.205
    40 = 43665.475
    41 = 108007.605
    42 = 83554.45
    43 = 96182.42
    44 = 95180.095
    45 = 109329.675
    46 = 114829.525
    47 = 88779.545
    48 = 38760.745
    49 = 142397.545
    50 = 56364.95
    51 = 67184.00499999999
    52 = 78621.545
    53 = 51128.56
    54 = 152362.6
    55 = 127560.75
    56 = 74122.73
    57 = 101667.67
    58 = 168114.745
    59 = 74192.55
    60 = 71275.02499999999
    61 = 131778.215
    62 = 93517.715
    63 = 49341.935
    64 = 47135.875
    65 = 125012.395
    66 = 64636.74000000001
    67 = 97223.795
    68 = 112624.38
    69 = 50734.125
    70 = 46520.64999999999
    71 = 96135.645
    72 = 86782.28999999999
    73 = 43128.655
    74 = 91917.825
    75 = 96481.045
    76 = 97003.55500000001
    77 = 127854.7
    78 = 82464.19499999999
    79 = 104810.445
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 2).Value = $values[$row]
}
